# Apply hybrid bold + color highlighting to quantitative metrics in the
# resume's PROFESSIONAL EXPERIENCE / KEY PROJECTS / KEY ACHIEVEMENTS
# bullet paragraphs, matching the target diff.

$d = $word.ActiveDocument

# Color used for the highlighted metrics: hex 2C3E50 -> BGR long for Font.Color
$HighlightColor = 5258796

function Apply-MetricHighlight {
    param(
        [string]$ParaText,
        [string[]]$BoldSegments
    )

    $paras = $d.Paragraphs
    $count = $paras.Count
    $matchIndex = -1

    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        # paragraph range text includes trailing paragraph mark (\r); compare
        # against that plus also allow exact match without it, to be safe.
        if ($t -eq ($ParaText + "`r") -or $t -eq $ParaText) {
            $matchIndex = $i
            break
        }
    }

    if ($matchIndex -eq -1) {
        Write-Output "WARNING: paragraph not found for: $ParaText"
        return
    }

    $p = $paras.Item($matchIndex)
    $pStart = $p.Range.Start

    $searchFrom = 0
    foreach ($seg in $BoldSegments) {
        $idx = $ParaText.IndexOf($seg, $searchFrom)
        if ($idx -lt 0) {
            Write-Output "WARNING: segment '$seg' not found in paragraph text"
            continue
        }
        $segStart = $pStart + $idx
        $segEnd = $segStart + $seg.Length
        $rng = $d.Range($segStart, $segEnd)
        $rng.Font.Bold = 1
        $rng.Font.Color = $HighlightColor
        $searchFrom = $idx + $seg.Length
    }
}

# 1. PROFESSIONAL EXPERIENCE - Siege Analytics - race coding errors bullet
Apply-MetricHighlight `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" `
    @("23%", "64%")

# 2. PROFESSIONAL EXPERIENCE - Siege Analytics - prediction accuracy bullet (long form)
Apply-MetricHighlight `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" `
    @("87%", "71%", "±4.2%", "±2.1%")

# 3. PROFESSIONAL EXPERIENCE - Senior Analyst - RFP vendors bullet
Apply-MetricHighlight `
    "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development" `
    @("1,200")

# 4. PROFESSIONAL EXPERIENCE - Programmer - meta-analysis framework bullet
Apply-MetricHighlight `
    "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" `
    @('$400M', '$1B')

# 5. KEY ACHIEVEMENTS - Algorithm mapping costs bullet
Apply-MetricHighlight `
    "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" `
    @("73.5%", '$4.7M')

# 6. KEY ACHIEVEMENTS - prediction accuracy bullet (short form)
Apply-MetricHighlight `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" `
    @("87%", "71%")

Write-Output "Done."
